# Auto-generated script to apply value updates to Malboro_Profits workbook
$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 2 (Leve Item ID 5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 131.35294
$ws.Range("I2").Value = 134.5
$ws.Range("J2").Value = 116.666664
$ws.Range("K2").Value = 134.5
$ws.Range("L2").Value = 116.666664
$ws.Range("M2").Value = -21.5
$ws.Range("N2").Value = -342.666664

# Hunk 1: sheet ALC, row 13 (Leve Item ID 2144)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2831

# Hunk 2: sheet ALC, row 18 (Leve Item ID 5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1649.8
$ws.Range("I18").Value = 1649.8
$ws.Range("K18").Value = 1649.8
$ws.Range("M18").Value = -1365.8

# Hunk 3: sheet ALC, row 32 (Leve Item ID 5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4799.5
$ws.Range("J32").Value = 4799.5
$ws.Range("L32").Value = 4799.5
$ws.Range("N32").Value = -5451.5

# Hunk 4: sheet ALC, row 33 (Leve Item ID 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 27787506
$ws.Range("I33").Value = 50001508
$ws.Range("J33").Value = 19999.75
$ws.Range("K33").Value = 50001508
$ws.Range("L33").Value = 19999.75
$ws.Range("M33").Value = -50001279
$ws.Range("N33").Value = -20457.75

# Hunk 5: sheet ALC, row 40 (Leve Item ID 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5736.625
$ws.Range("I40").Value = 3500
$ws.Range("J40").Value = 6056.143
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 6056.143
$ws.Range("M40").Value = -3325
$ws.Range("N40").Value = -6406.143

# Hunk 6: sheet ALC, row 43 (Leve Item ID 5472)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 10000
$ws.Range("K43").Value = 10000
$ws.Range("M43").Value = -9931

# Hunk 7: sheet ALC, row 99 (Leve Item ID 19883)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1562.5714
$ws.Range("I99").Value = 2114.25
$ws.Range("K99").Value = 6342.75
$ws.Range("M99").Value = -4844.75

# Hunk 8: sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 16524.55
$ws.Range("I132").Value = 14390.733
$ws.Range("K132").Value = 43172.199
$ws.Range("M132").Value = -40642.199

# Hunk 9: sheet ARM, row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4379.7646
$ws.Range("I2").Value = 920.125
$ws.Range("K2").Value = 920.125
$ws.Range("M2").Value = -807.125

# Hunk 10: sheet ARM, row 38 (Leve Item ID 2260)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 1929.6666
$ws.Range("I38").Value = 1929.6666
$ws.Range("K38").Value = 1929.6666
$ws.Range("M38").Value = -1462.6666

# Hunk 11: sheet ARM, row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14417.444
$ws.Range("I74").Value = 2712.1667
$ws.Range("K74").Value = 2712.1667
$ws.Range("M74").Value = -1838.1667

# Hunk 12: sheet ARM, row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14417.444
$ws.Range("I77").Value = 2712.1667
$ws.Range("K77").Value = 13560.8335
$ws.Range("M77").Value = -9192.833500000001

# Hunk 13: sheet ARM, row 102 (Leve Item ID 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 13552.474
$ws.Range("I102").Value = 2839.1333
$ws.Range("K102").Value = 2839.1333
$ws.Range("M102").Value = -1217.1333

# Hunk 14: sheet ARM, row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4379.7646
$ws.Range("I116").Value = 920.125
$ws.Range("K116").Value = 920.125
$ws.Range("M116").Value = 1373.875

# Hunk 15: sheet ARM, row 131 (Leve Item ID 34706)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080

# Hunk 16: sheet BSM, row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4379.7646
$ws.Range("I3").Value = 920.125
$ws.Range("K3").Value = 920.125
$ws.Range("M3").Value = -806.125

# Hunk 17: sheet BSM, row 107 (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2106.2666
$ws.Range("I107").Value = 1888.4445
$ws.Range("K107").Value = 1888.4445
$ws.Range("M107").Value = 31.55549999999994

# Hunk 18: sheet CRP, row 10 (Leve Item ID 1997)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 200168.8
$ws.Range("J10").Value = 500099.5
$ws.Range("L10").Value = 500099.5
$ws.Range("N10").Value = -500377.5

# Hunk 19: sheet CRP, row 22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3348.8
$ws.Range("I22").Value = 2873.5
$ws.Range("J22").Value = 3665.6667
$ws.Range("K22").Value = 2873.5
$ws.Range("L22").Value = 3665.6667
$ws.Range("M22").Value = -2523.5
$ws.Range("N22").Value = -4365.6667

# Hunk 20: sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 51685.363
$ws.Range("I31").Value = 27083
$ws.Range("K31").Value = 27083
$ws.Range("M31").Value = -26788

# Hunk 21: sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 51685.363
$ws.Range("I34").Value = 27083
$ws.Range("K34").Value = 27083
$ws.Range("M34").Value = -26881

# Hunk 22: sheet CRP, row 94 (Leve Item ID 32934)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 60006
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -60902

# Hunk 23: sheet CRP, row 107 (Leve Item ID 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3689.1365
$ws.Range("J107").Value = 5929
$ws.Range("L107").Value = 5929
$ws.Range("N107").Value = -9769

# Hunk 24: sheet CUL, row 33 (Leve Item ID 4867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 431.77777
$ws.Range("J33").Value = 569
$ws.Range("L33").Value = 3414
$ws.Range("N33").Value = -3980

# Hunk 25: sheet CUL, row 55 (Leve Item ID 4733)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3250
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3250
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9750
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -10104

# Hunk 26: sheet CUL, row 125 (Leve Item ID 36043)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5999
$ws.Range("I125").Value = 5999
$ws.Range("K125").Value = 17997
$ws.Range("M125").Value = -13077

# Hunk 27: sheet CUL, row 130 (Leve Item ID 36058)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 16124.25
$ws.Range("I130").Value = 3000
$ws.Range("K130").Value = 9000
$ws.Range("M130").Value = -3980

# Hunk 28: sheet CUL, row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1487.12
$ws.Range("J131").Value = 1499.3158
$ws.Range("L131").Value = 4497.9474
$ws.Range("N131").Value = -14577.9474

# Hunk 29: sheet GSM, row 34 (Leve Item ID 10924)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 49999
$ws.Range("J34").Value = 49999
$ws.Range("L34").Value = 49999
$ws.Range("N34").Value = -50535

# Hunk 30: sheet GSM, row 76 (Leve Item ID 10924)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 49999
$ws.Range("J76").Value = 49999
$ws.Range("L76").Value = 49999
$ws.Range("N76").Value = -50629

# Hunk 31: sheet GSM, row 79 (Leve Item ID 10924)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 49999
$ws.Range("J79").Value = 49999
$ws.Range("L79").Value = 49999
$ws.Range("N79").Value = -52183

# Hunk 32: sheet GSM, row 97 (Leve Item ID 19940)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 6116.6
$ws.Range("J97").Value = 14483.714
$ws.Range("L97").Value = 14483.714
$ws.Range("N97").Value = -15475.714

# Hunk 33: sheet GSM, row 100 (Leve Item ID 18367)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 47499
$ws.Range("J100").Value = 47499
$ws.Range("L100").Value = 47499
$ws.Range("N100").Value = -49663

# Hunk 34: sheet GSM, row 107 (Leve Item ID 27802)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 404
$ws.Range("I107").Value = 174.75
$ws.Range("J107").Value = 862.5
$ws.Range("K107").Value = 174.75
$ws.Range("L107").Value = 862.5
$ws.Range("M107").Value = 1745.25
$ws.Range("N107").Value = -4702.5

# Hunk 35: sheet GSM, row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 12134.869
$ws.Range("I126").Value = 14589.667
$ws.Range("K126").Value = 43769.001
$ws.Range("M126").Value = -41299.001

# Hunk 36: sheet GSM, row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13550.167
$ws.Range("I132").Value = 8211.429
$ws.Range("K132").Value = 24634.287
$ws.Range("M132").Value = -22104.287

# Hunk 37: sheet LTW, row 16 (Leve Item ID 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3619.7058
$ws.Range("I16").Value = 3733.4375
$ws.Range("K16").Value = 3733.4375
$ws.Range("M16").Value = -3563.4375

# Hunk 38: sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6111.0967
$ws.Range("I40").Value = 2912.125
$ws.Range("K40").Value = 2912.125
$ws.Range("M40").Value = -2776.125

# Hunk 39: sheet LTW, row 70 (Leve Item ID 10811)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 26160
$ws.Range("J70").Value = 26160
$ws.Range("L70").Value = 26160
$ws.Range("N70").Value = -26700

# Hunk 40: sheet LTW, row 73 (Leve Item ID 10811)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 26160
$ws.Range("J73").Value = 26160
$ws.Range("L73").Value = 26160
$ws.Range("N73").Value = -28032

# Hunk 41: sheet LTW, row 100 (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7496.0527
$ws.Range("I100").Value = 5241.7334
$ws.Range("K100").Value = 5241.7334
$ws.Range("M100").Value = -4700.7334

# Hunk 42: sheet LTW, row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7372.1333
$ws.Range("I122").Value = 5258
$ws.Range("J122").Value = 10543.333
$ws.Range("K122").Value = 15774
$ws.Range("L122").Value = 31629.999
$ws.Range("M122").Value = -13324
$ws.Range("N122").Value = -36529.999

# Hunk 43: sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2017220.8
$ws.Range("I132").Value = 9540.5
$ws.Range("K132").Value = 28621.5
$ws.Range("M132").Value = -26091.5

# Hunk 44: sheet LTW, row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 14696.027
$ws.Range("I136").Value = 17345.5
$ws.Range("K136").Value = 52036.5
$ws.Range("M136").Value = -49486.5

# Hunk 45: sheet WVR, row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 28876.625
$ws.Range("I126").Value = 15501
$ws.Range("J126").Value = 42252.25
$ws.Range("K126").Value = 46503
$ws.Range("L126").Value = 126756.75
$ws.Range("M126").Value = -44033
$ws.Range("N126").Value = -131696.75

# Hunk 46: sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9246.3125
$ws.Range("I132").Value = 3546.4211
$ws.Range("K132").Value = 10639.2633
$ws.Range("M132").Value = -8109.263300000001
